$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 539, shifting existing rows 539-567 down to 540-568.
$ws.Rows.Item(539).Insert()

# Populate the newly inserted row 539 with the new weekly data point.
$ws.Cells.Item(539, 1).Value = 9
$ws.Cells.Item(539, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(539, 3).Value = "Metropolitana"
$ws.Cells.Item(539, 4).Value = 45041
$ws.Cells.Item(539, 5).Value = 13
$ws.Cells.Item(539, 6).Value = 100112039
$ws.Cells.Item(539, 7).Value = "Ciboulette"
$ws.Cells.Item(539, 8).Value = "Sin especificar"
$ws.Cells.Item(539, 9).Value = "Primera"
$ws.Cells.Item(539, 10).Value = 340
$ws.Cells.Item(539, 11).Value = 1000
$ws.Cells.Item(539, 12).Value = 1200
$ws.Cells.Item(539, 13).Value = 1100
$ws.Cells.Item(539, 14).Value = '$/docena de atados'
$ws.Cells.Item(539, 15).Value = "Región Metropolitana"
$ws.Cells.Item(539, 16).Value = 367
$ws.Cells.Item(539, 17).Value = 3
$ws.Cells.Item(539, 18).Value = "Hortaliza"
